$wb = $excel.ActiveWorkbook

# --- Sheet 1: "Daily Orders" ---
$ws = $wb.Worksheets.Item("Daily Orders")

# Insert a new row above the current row 2 (pushes existing order rows down by one)
$ws.Rows.Item(2).Insert()

# Populate the newly inserted row 2 with the new order's data
$ws.Range("A2").Value = 5
$ws.Range("B2").Value = "2026-01-13 16:54"
$ws.Range("C2").Value = "Pooja"
$ws.Range("D2").Value = "'9096648553"
$ws.Range("E2").Value = "Level 1, Tower S3, CyberCity,`nMagarpatta City, Hadapsar, PUNE 411013"
$ws.Range("F2").Value = "Girl Holding Hands Thali x1"
$ws.Range("G2").Value = 0
$ws.Range("H2").Value = "NEW"
$ws.Range("I2").Value = "PENDING"
$ws.Range("J2").Value = "'"
$ws.Range("K2").Value = "'"
$ws.Range("L2").Value = "'"

# --- Sheet 2: "Summary" ---
$ws2 = $wb.Worksheets.Item("Summary")
$ws2.Range("A2").Value = 5
$ws2.Range("B2").Value = 4
